$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A67:A80").NumberFormat = "@"
$ws.Range("A67").Value = "2026-01-28"
$ws.Range("B67").Value = "15:33:56"
$ws.Range("C67").Value = "15:00"
$ws.Range("D67").Value = "Bathroom"
$ws.Range("E67").Value = "No Motion"
$ws.Range("F67").Value = "Inactive"
$ws.Range("A68").Value = "2026-01-28"
$ws.Range("B68").Value = "15:33:57"
$ws.Range("C68").Value = "15:00"
$ws.Range("D68").Value = "Bathroom"
$ws.Range("E68").Value = "No Motion"
$ws.Range("F68").Value = "Inactive"
$ws.Range("A69").Value = "2026-01-28"
$ws.Range("B69").Value = "15:33:59"
$ws.Range("C69").Value = "15:00"
$ws.Range("D69").Value = "Bathroom"
$ws.Range("E69").Value = "No Motion"
$ws.Range("F69").Value = "Inactive"
$ws.Range("A70").Value = "2026-01-28"
$ws.Range("B70").Value = "15:34:05"
$ws.Range("C70").Value = "15:00"
$ws.Range("D70").Value = "Bathroom"
$ws.Range("E70").Value = "No Motion"
$ws.Range("F70").Value = "Inactive"
$ws.Range("A71").Value = "2026-01-28"
$ws.Range("B71").Value = "15:34:10"
$ws.Range("C71").Value = "15:00"
$ws.Range("D71").Value = "Bathroom"
$ws.Range("E71").Value = "No Motion"
$ws.Range("F71").Value = "Inactive"
$ws.Range("A72").Value = "2026-01-28"
$ws.Range("B72").Value = "15:34:15"
$ws.Range("C72").Value = "15:00"
$ws.Range("D72").Value = "Bathroom"
$ws.Range("E72").Value = "No Motion"
$ws.Range("F72").Value = "Inactive"
$ws.Range("A73").Value = "2026-01-28"
$ws.Range("B73").Value = "15:34:20"
$ws.Range("C73").Value = "15:00"
$ws.Range("D73").Value = "Bathroom"
$ws.Range("E73").Value = "No Motion"
$ws.Range("F73").Value = "Inactive"
$ws.Range("A74").Value = "2026-01-28"
$ws.Range("B74").Value = "15:34:25"
$ws.Range("C74").Value = "15:00"
$ws.Range("D74").Value = "Bathroom"
$ws.Range("E74").Value = "No Motion"
$ws.Range("F74").Value = "Inactive"
$ws.Range("A75").Value = "2026-01-28"
$ws.Range("B75").Value = "15:34:30"
$ws.Range("C75").Value = "15:00"
$ws.Range("D75").Value = "Bathroom"
$ws.Range("E75").Value = "No Motion"
$ws.Range("F75").Value = "Inactive"
$ws.Range("A76").Value = "2026-01-28"
$ws.Range("B76").Value = "15:34:35"
$ws.Range("C76").Value = "15:00"
$ws.Range("D76").Value = "Bathroom"
$ws.Range("E76").Value = "No Motion"
$ws.Range("F76").Value = "Inactive"
$ws.Range("A77").Value = "2026-01-28"
$ws.Range("B77").Value = "15:34:40"
$ws.Range("C77").Value = "15:00"
$ws.Range("D77").Value = "Bathroom"
$ws.Range("E77").Value = "No Motion"
$ws.Range("F77").Value = "Inactive"
$ws.Range("A78").Value = "2026-01-28"
$ws.Range("B78").Value = "15:34:45"
$ws.Range("C78").Value = "15:00"
$ws.Range("D78").Value = "Bathroom"
$ws.Range("E78").Value = "No Motion"
$ws.Range("F78").Value = "Inactive"
$ws.Range("A79").Value = "2026-01-28"
$ws.Range("B79").Value = "15:34:50"
$ws.Range("C79").Value = "15:00"
$ws.Range("D79").Value = "Bathroom"
$ws.Range("E79").Value = "No Motion"
$ws.Range("F79").Value = "Inactive"
$ws.Range("A80").Value = "2026-01-28"
$ws.Range("B80").Value = "15:34:55"
$ws.Range("C80").Value = "15:00"
$ws.Range("D80").Value = "Bathroom"
$ws.Range("E80").Value = "No Motion"
$ws.Range("F80").Value = "Inactive"
$ws.Range("A67:A80").NumberFormat = "@"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A68:A82").NumberFormat = "@"
$ws.Range("E68:E82").NumberFormat = "@"
$ws.Range("A68").Value = "2026-01-28"
$ws.Range("B68").Value = "15:33:56"
$ws.Range("C68").Value = "15:00"
$ws.Range("D68").Value = "Bathroom"
$ws.Range("E68").Value = "88.4%"
$ws.Range("F68").Value = "Active"
$ws.Range("A69").Value = "2026-01-28"
$ws.Range("B69").Value = "15:33:56"
$ws.Range("C69").Value = "15:00"
$ws.Range("D69").Value = "Bathroom"
$ws.Range("E69").Value = "87.5%"
$ws.Range("F69").Value = "Active"
$ws.Range("A70").Value = "2026-01-28"
$ws.Range("B70").Value = "15:33:58"
$ws.Range("C70").Value = "15:00"
$ws.Range("D70").Value = "Bathroom"
$ws.Range("E70").Value = "88.4%"
$ws.Range("F70").Value = "Active"
$ws.Range("A71").Value = "2026-01-28"
$ws.Range("B71").Value = "15:34:02"
$ws.Range("C71").Value = "15:00"
$ws.Range("D71").Value = "Bathroom"
$ws.Range("E71").Value = "87.5%"
$ws.Range("F71").Value = "Active"
$ws.Range("A72").Value = "2026-01-28"
$ws.Range("B72").Value = "15:34:06"
$ws.Range("C72").Value = "15:00"
$ws.Range("D72").Value = "Bathroom"
$ws.Range("E72").Value = "88.4%"
$ws.Range("F72").Value = "Active"
$ws.Range("A73").Value = "2026-01-28"
$ws.Range("B73").Value = "15:34:10"
$ws.Range("C73").Value = "15:00"
$ws.Range("D73").Value = "Bathroom"
$ws.Range("E73").Value = "88.4%"
$ws.Range("F73").Value = "Active"
$ws.Range("A74").Value = "2026-01-28"
$ws.Range("B74").Value = "15:34:14"
$ws.Range("C74").Value = "15:00"
$ws.Range("D74").Value = "Bathroom"
$ws.Range("E74").Value = "87.5%"
$ws.Range("F74").Value = "Active"
$ws.Range("A75").Value = "2026-01-28"
$ws.Range("B75").Value = "15:34:18"
$ws.Range("C75").Value = "15:00"
$ws.Range("D75").Value = "Bathroom"
$ws.Range("E75").Value = "88.4%"
$ws.Range("F75").Value = "Active"
$ws.Range("A76").Value = "2026-01-28"
$ws.Range("B76").Value = "15:34:26"
$ws.Range("C76").Value = "15:00"
$ws.Range("D76").Value = "Bathroom"
$ws.Range("E76").Value = "88.5%"
$ws.Range("F76").Value = "Active"
$ws.Range("A77").Value = "2026-01-28"
$ws.Range("B77").Value = "15:34:30"
$ws.Range("C77").Value = "15:00"
$ws.Range("D77").Value = "Bathroom"
$ws.Range("E77").Value = "88.5%"
$ws.Range("F77").Value = "Active"
$ws.Range("A78").Value = "2026-01-28"
$ws.Range("B78").Value = "15:34:34"
$ws.Range("C78").Value = "15:00"
$ws.Range("D78").Value = "Bathroom"
$ws.Range("E78").Value = "87.6%"
$ws.Range("F78").Value = "Active"
$ws.Range("A79").Value = "2026-01-28"
$ws.Range("B79").Value = "15:34:38"
$ws.Range("C79").Value = "15:00"
$ws.Range("D79").Value = "Bathroom"
$ws.Range("E79").Value = "88.5%"
$ws.Range("F79").Value = "Active"
$ws.Range("A80").Value = "2026-01-28"
$ws.Range("B80").Value = "15:34:46"
$ws.Range("C80").Value = "15:00"
$ws.Range("D80").Value = "Bathroom"
$ws.Range("E80").Value = "87.6%"
$ws.Range("F80").Value = "Active"
$ws.Range("A81").Value = "2026-01-28"
$ws.Range("B81").Value = "15:34:50"
$ws.Range("C81").Value = "15:00"
$ws.Range("D81").Value = "Bathroom"
$ws.Range("E81").Value = "88.6%"
$ws.Range("F81").Value = "Active"
$ws.Range("A82").Value = "2026-01-28"
$ws.Range("B82").Value = "15:34:54"
$ws.Range("C82").Value = "15:00"
$ws.Range("D82").Value = "Bathroom"
$ws.Range("E82").Value = "87.6%"
$ws.Range("F82").Value = "Active"
$ws.Range("A68:A82").NumberFormat = "@"
$ws.Range("E68:E82").NumberFormat = "@"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A68:A82").NumberFormat = "@"
$ws.Range("A68").Value = "2026-01-28"
$ws.Range("B68").Value = "15:33:56"
$ws.Range("C68").Value = "15:00"
$ws.Range("D68").Value = "Bathroom"
$ws.Range("E68").Value = "22.9C"
$ws.Range("F68").Value = "Active"
$ws.Range("A69").Value = "2026-01-28"
$ws.Range("B69").Value = "15:33:57"
$ws.Range("C69").Value = "15:00"
$ws.Range("D69").Value = "Bathroom"
$ws.Range("E69").Value = "22.9C"
$ws.Range("F69").Value = "Active"
$ws.Range("A70").Value = "2026-01-28"
$ws.Range("B70").Value = "15:33:58"
$ws.Range("C70").Value = "15:00"
$ws.Range("D70").Value = "Bathroom"
$ws.Range("E70").Value = "22.9C"
$ws.Range("F70").Value = "Active"
$ws.Range("A71").Value = "2026-01-28"
$ws.Range("B71").Value = "15:34:02"
$ws.Range("C71").Value = "15:00"
$ws.Range("D71").Value = "Bathroom"
$ws.Range("E71").Value = "22.9C"
$ws.Range("F71").Value = "Active"
$ws.Range("A72").Value = "2026-01-28"
$ws.Range("B72").Value = "15:34:06"
$ws.Range("C72").Value = "15:00"
$ws.Range("D72").Value = "Bathroom"
$ws.Range("E72").Value = "22.9C"
$ws.Range("F72").Value = "Active"
$ws.Range("A73").Value = "2026-01-28"
$ws.Range("B73").Value = "15:34:10"
$ws.Range("C73").Value = "15:00"
$ws.Range("D73").Value = "Bathroom"
$ws.Range("E73").Value = "22.9C"
$ws.Range("F73").Value = "Active"
$ws.Range("A74").Value = "2026-01-28"
$ws.Range("B74").Value = "15:34:14"
$ws.Range("C74").Value = "15:00"
$ws.Range("D74").Value = "Bathroom"
$ws.Range("E74").Value = "22.9C"
$ws.Range("F74").Value = "Active"
$ws.Range("A75").Value = "2026-01-28"
$ws.Range("B75").Value = "15:34:18"
$ws.Range("C75").Value = "15:00"
$ws.Range("D75").Value = "Bathroom"
$ws.Range("E75").Value = "22.9C"
$ws.Range("F75").Value = "Active"
$ws.Range("A76").Value = "2026-01-28"
$ws.Range("B76").Value = "15:34:26"
$ws.Range("C76").Value = "15:00"
$ws.Range("D76").Value = "Bathroom"
$ws.Range("E76").Value = "22.9C"
$ws.Range("F76").Value = "Active"
$ws.Range("A77").Value = "2026-01-28"
$ws.Range("B77").Value = "15:34:30"
$ws.Range("C77").Value = "15:00"
$ws.Range("D77").Value = "Bathroom"
$ws.Range("E77").Value = "22.9C"
$ws.Range("F77").Value = "Active"
$ws.Range("A78").Value = "2026-01-28"
$ws.Range("B78").Value = "15:34:35"
$ws.Range("C78").Value = "15:00"
$ws.Range("D78").Value = "Bathroom"
$ws.Range("E78").Value = "22.9C"
$ws.Range("F78").Value = "Active"
$ws.Range("A79").Value = "2026-01-28"
$ws.Range("B79").Value = "15:34:38"
$ws.Range("C79").Value = "15:00"
$ws.Range("D79").Value = "Bathroom"
$ws.Range("E79").Value = "22.9C"
$ws.Range("F79").Value = "Active"
$ws.Range("A80").Value = "2026-01-28"
$ws.Range("B80").Value = "15:34:47"
$ws.Range("C80").Value = "15:00"
$ws.Range("D80").Value = "Bathroom"
$ws.Range("E80").Value = "22.9C"
$ws.Range("F80").Value = "Active"
$ws.Range("A81").Value = "2026-01-28"
$ws.Range("B81").Value = "15:34:51"
$ws.Range("C81").Value = "15:00"
$ws.Range("D81").Value = "Bathroom"
$ws.Range("E81").Value = "22.9C"
$ws.Range("F81").Value = "Active"
$ws.Range("A82").Value = "2026-01-28"
$ws.Range("B82").Value = "15:34:55"
$ws.Range("C82").Value = "15:00"
$ws.Range("D82").Value = "Bathroom"
$ws.Range("E82").Value = "22.9C"
$ws.Range("F82").Value = "Active"
$ws.Range("A68:A82").NumberFormat = "@"

$ws = $wb.Worksheets.Item("ALERTS")
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-01-28"
$ws.Range("B2").Value = "15:34:50"
$ws.Range("C2").Value = "15:00"
$ws.Range("D2").Value = "Bathroom Door"
$ws.Range("E2").Value = "ENTER"
$ws.Range("F2").Value = "User ENTERED Bathroom"
$ws.Range("A2").NumberFormat = "@"

